$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216, pushing the existing rows 216-266 down to 217-267.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new record.
$ws.Range("A216").Value = 10
$ws.Range("B216").Value = "Vega Modelo de Temuco"
$ws.Range("C216").Value = "La Araucanía"
$ws.Range("D216").Value = 44855
$ws.Range("E216").Value = 9
$ws.Range("F216").Value = 100112013
$ws.Range("G216").Value = "Alcachofa"
$ws.Range("H216").Value = "Española"
$ws.Range("I216").Value = "Extra"
$ws.Range("J216").Value = 70
$ws.Range("K216").Value = 13000
$ws.Range("L216").Value = 14000
$ws.Range("M216").Value = 13571
$ws.Range("N216").Value = "`$/caja 30 unidades"
$ws.Range("O216").Value = "Región Metropolitana"
$ws.Range("P216").Value = 452
$ws.Range("Q216").Value = 30
$ws.Range("R216").Value = "Hortaliza"
